$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Message Chains" task with "Long Method (counting statements)" for Andrii
$ws.Range("C1").Value = "Long Method (counting statements)"

$ws.Range("C1").Select()
